$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 99, shifting existing rows (99..164) down to (100..165).
$ws.Rows.Item(99).Insert()

# Populate the newly inserted row 99 with the new daily price record.
$ws.Range("A99").Value = 9
$ws.Range("B99").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C99").Value = "Metropolitana"
$ws.Range("D99").Value = 45068
$ws.Range("E99").Value = 13
$ws.Range("F99").Value = "Fruta"
$ws.Range("G99").Value = 100101
$ws.Range("H99").Value = "Berries"
$ws.Range("I99").Value = 100101004
$ws.Range("J99").Value = "Frambuesa"
$ws.Range("K99").Value = "Sin especificar"
$ws.Range("L99").Value = "Primera"
$ws.Range("M99").Value = 390
$ws.Range("N99").Value = 8000
$ws.Range("O99").Value = 9000
$ws.Range("P99").Value = 8462
$ws.Range("Q99").Value = "$/bandeja 2 kilos"
$ws.Range("R99").Value = "Provincia de Linares"
$ws.Range("S99").Value = 4231
$ws.Range("T99").Value = 2
